$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete medway_booking_m (row 19), keep row 18 (medway_booking_y) and rename it to a wildcard pattern
$ws.Rows.Item(19).Delete()
$ws.Range("A18").Value = "BiBBS_CohortInfo.pregnancy.medway_booking_?"

# Rename fin_2cutmeals (row 14) to a wildcard pattern
$ws.Range("A14").Value = "BiBBS_Baseline.pregnancy_survey.fin_?cutmeals"

# Collapse the four mes1_* rows (5-8) into a single wildcard row
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
$ws.Range("A5").Value = "BiBBS_Baseline.pregnancy_survey.mes1_*"

[void]$ws.Range("B8").Select()
